$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix misconception of TTD:PDSpec relationship.
# Already satisfied by Equipment:PDSpec relationship, so mark both
# directions of the TTD <-> Pressure Drop Specification matrix cell as "None".

$ws.Range("S9").Value = "None"
$ws.Range("I19").Value = "None"
